$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-looking decimal numbers as text in
# the source workbook (e.g. "413.51", "1.00", "0.0000119"). Writing those
# strings straight into a General-formatted cell makes Excel silently
# reinterpret them as numbers and drop formatting (trailing zeros, etc.),
# so first mark the cells as Text, then restore the default "Normal" style
# afterwards so no stray number-format style is left on the cell.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.831.36'
$ws.Range("E2").Value = '  +8.65%  '
$ws.Range("D3").Value = '3.355.14'
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '413.51'
$ws.Range("E5").Value = '  +4.27%  '
$ws.Range("D6").Value = '116.06'
$ws.Range("E6").Value = '  +7.15%  '
$ws.Range("D7").Value = '3.348.23'
$ws.Range("E7").Value = '  +3.43%  '
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  +20.17%  '
$ws.Range("D12").Value = '39.87'
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '3.871.97'
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").Value = '8.31'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '19.26'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("D17").Value = '3.340.91'
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").Value = '61.477.78'
$ws.Range("E18").Value = '  +8.23%  '
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = '10.81'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '0.0000119'
$ws.Range("E21").Value = '  +9.26%  '
$ws.Range("D22").Value = '3.33'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '12.51'
$ws.Range("E23").Value = '  -4.25%  '
$ws.Range("D24").Value = '294.77'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '74.59'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("D27").Value = '29.30'
$ws.Range("E27").Value = '  +4.22%  '
$ws.Range("D28").Value = '7.92'
$ws.Range("E28").Value = '  +9.32%  '
$ws.Range("D29").Value = '0.174'
$ws.Range("E29").Value = '  +3.43%  '
$ws.Range("D30").Value = '4.26'
$ws.Range("E30").Value = '  -1.97%  '
$ws.Range("D31").Value = '7.62'
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '42.88'
$ws.Range("E32").Value = '  +7.31%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.113'
$ws.Range("E33").Value = '  +4.64%  '
$ws.Range("D34").Value = '11.44'
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").Value = '2.54'
$ws.Range("E35").Value = '  +19.13%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '0.0491'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = '52.24'
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("D39").Value = '3.13'
$ws.Range("E39").Value = '  +6.60%  '
$ws.Range("D40").Value = '0.996'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("D42").Value = '133.24'
$ws.Range("E42").Value = '  -4.80%  '
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '0.286'
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").Value = '3.87'
$ws.Range("E46").Value = '  -2.11%  '
$ws.Range("D47").Value = '16.50'
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").Value = '2.23'
$ws.Range("E48").Value = '  +2.44%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.168.52'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '21.14'
$ws.Range("E50").Value = '  -4.92%  '
$ws.Range("D51").Value = '3.679.23'
$ws.Range("E51").Value = '  +2.94%  '

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
